$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.176486372947693
$ws.Range("B1").Value = 2.414991140365601
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.337298154830933
$ws.Range("E1").Value = 1.204263210296631
